$d = $word.ActiveDocument

# Step 1: Center-align the (only) paragraph in the document.
# (wdAlignParagraphCenter = 1)
$d.Paragraphs(1).Alignment = 1

# Step 2: "the idea of vr was first" -> "the idea of VR technology was first"
$rng = $d.Content
$ok = $rng.Find.Execute("the idea of vr was first", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $ok) { throw "Step 2 Find failed" }
$rng.Text = "the idea of VR technology was first"

# Step 3: "Stanley G.Weinbaum in" -> "Stanley G. Weinbaum in"
$rng = $d.Content
$ok = $rng.Find.Execute("Stanley G.Weinbaum in", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $ok) { throw "Step 3 Find failed" }
$rng.Text = "Stanley G. Weinbaum in"

# Step 4: Rewrite the closing sentence and append the new paragraphs of
# content. This range stops right before the existing "_GoBack" bookmark,
# which therefore stays untouched and keeps its original position (it will
# end up sitting right at the end of the paragraph once step 5 below
# clears out the old trailing run).
$rng = $d.Content
$ok = $rng.Find.Execute("patented a design in 1960 called The Sword of Damocles and it is considered the first virtual reality head-mounted display. The drawings of it look incredibly like", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $ok) { throw "Step 4 Find failed" }
$rng.Text = "patented a design in 1960 called The Sword of Damocles and it is considered the first virtual reality head-mounted display. The drawings of it look incredibly like the VR technology we see today. As the years progressed a man named Ivan Sutherland detailed the idea of an “Ultimate Display” in 1965. This device would be able to simulate a virtual reality to the point that the user would not be able to tell the difference from actual reality. “The ultimate display would, of course, be a room within which the computer can control the existence of matter. A chair displayed in such a room would be good enough to sit in. Handcuffs displayed in such a room would be confining, and a bullet displayed in such a room would be fatal. With appropriate programming such a display could literally be the Wonderland into which Alice walked.” – Ivan Sutherland. Obviously as it was only 1965, the technology available made it impossible to create such an impressive piece of technology. The idea of a virtual world indistinguishable from our own, sparked ideas in the minds of inventors in the years to come. Like American computer artist Myron W. Krueger, who created Videoplace. This was an artificial reality that surrounded the user, it also responded to their movements without needing goggles or gloves, the user would be projected onto a screen in front of them, they could change the image of themselves by moving around. "

# Step 5: The old trailing text (previously " the vr technology we see
# today.", now superseded by the rewritten text inserted in step 4) is
# removed so nothing remains after the bookmark. MatchCase is used so this
# cannot accidentally match the (differently-cased) "VR" text introduced
# above.
$rng = $d.Content
$ok = $rng.Find.Execute(" the vr technology we see today.", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $ok) { throw "Step 5 Find failed" }
$rng.Text = ""

Write-Host $d.Content.Text
